$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.570639666666668
$ws.Range("H2").Value = 16.711919
$ws.Range("I2").Value = 0.06371389580191485
$ws.Range("J2").Value = 0.06371389580191483
$ws.Range("M2").Value = 9.101794333333332
$ws.Range("N2").Value = 27.305383
$ws.Range("O2").Value = 0.1526015110517656
$ws.Range("P2").Value = 0.1526015110517656
$ws.Range("Q2").Value = 50.70281655110856
$ws.Range("R2").Value = 456.325348959977
$ws.Range("S2").Value = 0.009722836774366949
$ws.Range("T2").Value = 0.009722836774366948
$ws.Range("G3").Value = 5.570639666666668
$ws.Range("H3").Value = 16.711919
$ws.Range("I3").Value = 0.06371389580191485
$ws.Range("J3").Value = 0.06371389580191483
$ws.Range("O3").Value = 0.5991759712230392
$ws.Range("P3").Value = 0.5991759712230392
$ws.Range("Q3").Value = 199.0800034768238
$ws.Range("R3").Value = 1791.720031291414
$ws.Range("S3").Value = 0.03817583539751585
$ws.Range("T3").Value = 0.03817583539751584
$ws.Range("G4").Value = 5.570639666666668
$ws.Range("H4").Value = 16.711919
$ws.Range("I4").Value = 0.06371389580191485
$ws.Range("J4").Value = 0.06371389580191483
$ws.Range("O4").Value = 0.2482225177251951
$ws.Range("P4").Value = 0.2482225177251951
$ws.Range("Q4").Value = 82.47350038234924
$ws.Range("R4").Value = 742.2615034411431
$ws.Range("S4").Value = 0.01581522363003204
$ws.Range("T4").Value = 0.01581522363003204
$ws.Range("I5").Value = 0.5612617414829089
$ws.Range("J5").Value = 0.5612617414829089
$ws.Range("M5").Value = 9.101794333333332
$ws.Range("N5").Value = 27.305383
$ws.Range("O5").Value = 0.1526015110517656
$ws.Range("P5").Value = 0.1526015110517656
$ws.Range("Q5").Value = 446.6459122832101
$ws.Range("R5").Value = 4019.813210548891
$ws.Range("S5").Value = 0.08564938984583732
$ws.Range("T5").Value = 0.08564938984583732
$ws.Range("I6").Value = 0.5612617414829089
$ws.Range("J6").Value = 0.5612617414829089
$ws.Range("O6").Value = 0.5991759712230392
$ws.Range("P6").Value = 0.5991759712230392
$ws.Range("S6").Value = 0.3362945490633563
$ws.Range("T6").Value = 0.3362945490633563
$ws.Range("I7").Value = 0.5612617414829089
$ws.Range("J7").Value = 0.5612617414829089
$ws.Range("O7").Value = 0.2482225177251951
$ws.Range("P7").Value = 0.2482225177251951
$ws.Range("S7").Value = 0.1393178025737152
$ws.Range("T7").Value = 0.1393178025737152
$ws.Range("I8").Value = 0.3750243627151764
$ws.Range("J8").Value = 0.3750243627151764
$ws.Range("M8").Value = 9.101794333333332
$ws.Range("N8").Value = 27.305383
$ws.Range("O8").Value = 0.1526015110517656
$ws.Range("P8").Value = 0.1526015110517656
$ws.Range("Q8").Value = 298.4402574292517
$ws.Range("R8").Value = 2685.962316863266
$ws.Range("S8").Value = 0.05722928443156134
$ws.Range("T8").Value = 0.05722928443156134
$ws.Range("I9").Value = 0.3750243627151764
$ws.Range("J9").Value = 0.3750243627151764
$ws.Range("O9").Value = 0.5991759712230392
$ws.Range("P9").Value = 0.5991759712230392
$ws.Range("S9").Value = 0.2247055867621672
$ws.Range("T9").Value = 0.2247055867621672
$ws.Range("I10").Value = 0.3750243627151764
$ws.Range("J10").Value = 0.3750243627151764
$ws.Range("O10").Value = 0.2482225177251951
$ws.Range("P10").Value = 0.2482225177251951
$ws.Range("S10").Value = 0.09308949152144788
$ws.Range("T10").Value = 0.09308949152144788
